# MIUA 2018 data refresh: pre/post-length source values recalculated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 22
$ws.Range("G3").Value = 59
$ws.Range("G4").Value = 63
$ws.Range("G5").Value = 63
$ws.Range("D6").Value = 141
$ws.Range("G6").Value = 168
$ws.Range("D7").Value = 141
$ws.Range("D10").Value = 12
$ws.Range("G10").Value = 39
$ws.Range("D11").Value = 12
$ws.Range("G11").Value = 39
$ws.Range("D12").Value = 32
$ws.Range("G12").Value = 60
$ws.Range("D13").Value = 32
$ws.Range("G13").Value = 60
$ws.Range("G14").Value = 67
$ws.Range("D15").Value = 298
$ws.Range("G15").Value = 322
$ws.Range("D16").Value = 298
$ws.Range("G16").Value = 322
$ws.Range("D19").Value = 146
$ws.Range("D21").Value = 207
$ws.Range("G21").Value = 249
$ws.Range("D22").Value = 207
$ws.Range("G22").Value = 249
$ws.Range("G24").Value = 390
$ws.Range("G25").Value = 50
$ws.Range("G26").Value = 50
$ws.Range("D27").Value = 45
$ws.Range("G27").Value = 78
$ws.Range("D28").Value = 45
$ws.Range("G28").Value = 78
$ws.Range("D29").Value = 73
$ws.Range("G29").Value = 106
$ws.Range("D30").Value = 73
$ws.Range("D31").Value = 111
$ws.Range("D32").Value = 111
$ws.Range("G33").Value = 143
$ws.Range("G35").Value = 245
$ws.Range("G36").Value = 245
$ws.Range("G37").Value = 57
$ws.Range("D38").Value = 388
$ws.Range("G38").Value = 414

# View refresh to match the re-saved worksheet (zoom + active selection).
$ws.Activate()
$excel.ActiveWindow.Zoom = 110
$ws.Range("J10").Select()
